$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# Insert 6 fresh rows above the old row 43 ("DenseDepth" block): 5 rows for the
# new "DORN (Weighted histogram matching)" ablation block + 1 blank separator
# row, mirroring the blank row that already separates row 41 from the old 43.
$ws.Range("A43:A48").EntireRow.Insert()

# --- New block: DORN (Weighted histogram matching) ablation, rows 43-47 ---
$ws.Range("A43").Value = "DORN (Weighted histogram matching)"
$ws.Range("B43").Value = 0.9045
$ws.Range("C43").Value = 0.9705
$ws.Range("D43").Value = 0.9892
$ws.Range("F43").Value = 0.4142
$ws.Range("G43").Value = 0.0912
$ws.Range("I43").Value = 0.0395
$ws.Range("L43").Value = "Intensity Only"

$ws.Range("A44").Value = "DORN (Weighted histogram matching)"
$ws.Range("B44").Value = 0.9045
$ws.Range("C44").Value = 0.9705
$ws.Range("D44").Value = 0.9892
$ws.Range("F44").Value = 0.4143
$ws.Range("G44").Value = 0.0912
$ws.Range("I44").Value = 0.0395
$ws.Range("L44").Value = "Intensity and Falloff"

$ws.Range("A45").Value = "DORN (Weighted histogram matching)"
$ws.Range("B45").Value = 0.9041
$ws.Range("C45").Value = 0.97
$ws.Range("D45").Value = 0.9889
$ws.Range("F45").Value = 0.4173
$ws.Range("G45").Value = 0.0902
$ws.Range("I45").Value = 0.0396
$ws.Range("L45").Value = "Intensity, Falloff, and DC/Ambient"

$ws.Range("A46").Value = "DORN (Weighted histogram matching)"
$ws.Range("B46").Value = 0.9041
$ws.Range("C46").Value = 0.9701
$ws.Range("D46").Value = 0.9889
$ws.Range("F46").Value = 0.4168
$ws.Range("G46").Value = 0.0903
$ws.Range("I46").Value = 0.0396
$ws.Range("L46").Value = "Intensity, Falloff, DC/Ambient, and Jitter"

$ws.Range("A47").Value = "DORN (Weighted histogram matching)"
$ws.Range("B47").Value = 0.9031
$ws.Range("C47").Value = 0.9691
$ws.Range("D47").Value = 0.9881
$ws.Range("F47").Value = 0.4586
$ws.Range("G47").Value = 0.0926
$ws.Range("I47").Value = 0.04
$ws.Range("L47").Value = "Intensity, Falloff, DC/Ambient, Jitter, and Poisson Noise"

# Row 48 stays blank (separator), same role as the blank row 42 above.

# Move the view: scroll so row 27 is at the top, and select the new trailing
# blank row (now row 52) the way the author's last save left it.
$ws.Range("A52:XFD52").Select()
$excel.ActiveWindow.ScrollRow = 27
